$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Sl. No"
$ws.Range("B1").Value = "Table name"

$ws.Range("B1").Select()
